$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("A2").Value = "bef653ef-1498-4e68-aacf-b9978b2697c6.md"
$wsOverview.Range("B2").Value = "e2e\bef653ef-1498-4e68-aacf-b9978b2697c6.md"
$wsOverview.Range("G2").Value = "2017-02-28 08:14:19"

$wsZhCn.Range("A2").Value = "bef653ef-1498-4e68-aacf-b9978b2697c6.md"
$wsZhCn.Range("G2").Value = "bef653ef-1498-4e68-aacf-b9978b2697c6.da04188e841ab7dc0fbddef15463dd9b2f08fcfb.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2017-02-28 08:14:04"

$wsDeDe.Range("A2").Value = "bef653ef-1498-4e68-aacf-b9978b2697c6.md"
$wsDeDe.Range("G2").Value = "bef653ef-1498-4e68-aacf-b9978b2697c6.da04188e841ab7dc0fbddef15463dd9b2f08fcfb.de-de.xlf"
$wsDeDe.Range("H2").Value = "2017-02-28 08:14:19"
